# Sara Alert invalid monitorees template:
# rename the "Group Number" column header to "Assigned User"
# and nudge the sheet's horizontal scroll position one column to the left
# (topLeftCell CE1 -> CD1), keeping the existing selection (CS2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared-string header text in the last column (CS1).
$ws.Range("CS1").Value = "Assigned User"

# Keep the selection on CS2 and scroll the window so column CD is the
# left-most visible column (was CE).
$ws.Range("CS2").Select()
$excel.ActiveWindow.ScrollColumn = 82
